$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.23%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.56%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.113"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.45%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08006"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.98%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.506"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.74%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.647"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.07%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.085"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "17.02%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1295"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "6.23%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1901"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.72%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09390"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.06%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04221"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.80%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1038"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.26%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.74%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005946"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.20%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.377"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.84%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3368"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.47%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.096"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.94%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.43%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3136"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.56%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04200"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.49%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001275"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.65%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004581"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "15.04%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001342"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.03%"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "10.33%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05398"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.76%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005634"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-11.88%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007745"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.64%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1414"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.47%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007342"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.91%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007876"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.49%"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.00%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006753"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.51%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000745"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.65%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05691"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "22.98%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003973"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.44%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002087"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.65%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001988"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.65%"
